$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell as TEXT (shared string),
# without Excel's auto number-detection turning it into a Number, and
# without minting a new cell style on the destination. We stage a formula
# that evaluates to the literal string on a scratch cell, copy it, and
# paste only the VALUE into the destination, so the destination keeps its
# own existing style untouched.
function Set-TextValue {
    param($range, [string]$text)
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = 0
}

# --- Update existing rows ---
Set-TextValue $ws.Range("A2") "10007973"
$ws.Range("B2").Value = "BUAVITA MANGGA 245ML"

Set-TextValue $ws.Range("A3") "10007970"
$ws.Range("B3").Value = "BUAVITA JAMBU SL 245"

Set-TextValue $ws.Range("A6") "20014069"
$ws.Range("B6").Value = "FF UHT FULL CRM 946"
$ws.Range("F6").Value = "RT,(E-1B)"

# --- Add new row 7 (same look as row 6: copy formatting first) ---
$ws.Range("A6:F6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Set-TextValue $ws.Range("A7") "20019674"
$ws.Range("B7").Value = "YOU C1000 ORG WTR500"
$ws.Range("C7").Value = "DU1AMKT"
Set-TextValue $ws.Range("D7") "1"
Set-TextValue $ws.Range("E7") "6"
$ws.Range("F7").Value = "RT,(E-1B)"

# Clean up the scratch cell / column so it doesn't linger in the sheet.
$ws.Range("Z100").EntireColumn.Delete()
